$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "53-43=",
    "76+5=",
    "19+49=",
    "58-53=",
    "3+37=",
    "11+54=",
    "97-96=",
    "33+19=",
    "2+86=",
    "87-16=",
    "22+12=",
    "49+15=",
    "95-77=",
    "16+4=",
    "46-25=",
    "31+1=",
    "60-4=",
    "55-17=",
    "37+2=",
    "25+41=",
    "60+18=",
    "67-4=",
    "79-27=",
    "80-42=",
    "48-14=",
    "78-8=",
    "93-88=",
    "64-42=",
    "54+44=",
    "30+37=",
    "89-18=",
    "87-22=",
    "99-39=",
    "94-81=",
    "12+72=",
    "2+35=",
    "67-3=",
    "69+24=",
    "76-19=",
    "8+14=",
    "11+47=",
    "62-21=",
    "50+24=",
    "45+37=",
    "72+2=",
    "58+15=",
    "13+83=",
    "24+56=",
    "25+32=",
    "91+8=",
    "29+28=",
    "64+6=",
    "88-44=",
    "23+74=",
    "81-5=",
    "53+38=",
    "17+80=",
    "19+27=",
    "39+17=",
    "7+42=",
    "43+36=",
    "47-9=",
    "66-48=",
    "80-60=",
    "6+72=",
    "41+31=",
    "61-49=",
    "43-10=",
    "15+55=",
    "36-34=",
    "63-43=",
    "53-49=",
    "34-27=",
    "54+16=",
    "95-7=",
    "4+17=",
    "2+47=",
    "46-15=",
    "45+33=",
    "60+9=",
    "22+53=",
    "9+72=",
    "47+38=",
    "63-48=",
    "43+44=",
    "84-16=",
    "68+30=",
    "38+12=",
    "34+65=",
    "71-66=",
    "68+24=",
    "14+54=",
    "85-69=",
    "90-6=",
    "59+14=",
    "15+6=",
    "16+81=",
    "59+4=",
    "7+1=",
    "23+18="
)

$cols = 5
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = [Math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $values[$i]
}
Write-Host "Done updating $($values.Count) cells"